$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the extent of the data (header row + all data rows/cols).
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# Rename header columns: "<name>_old" -> "<name>_FV2210", "<name>_new" -> "<name>_FV2304"
# (the workbook compares an old format version, "FV2210", against a new one, "FV2304")
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# Turn the data range into a real Excel Table (ListObject) named "Table1"
$range = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
